$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "2022" column (L) mirroring the existing 2021 column (K) ---
$ws.Range("K3:K8").Copy()
$ws.Range("L3:L8").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("L3").Value = 2022
$ws.Range("L4").Formula = "=L5/L6*1000"
$ws.Range("L5").Value = 12673.2
$ws.Range("L6").Value = 7037.6
$ws.Range("L7").Value = 51.3
$ws.Range("L8").Value = 0.1

# The whole new data column is right aligned, and the computed "per capita" row is bold
$ws.Range("L4:L8").HorizontalAlignment = -4152   # xlRight
$ws.Range("L4").Font.Bold = $true

# --- The 2014 figures for the last two indicators were never collected ---
$ws.Range("D7").Value = "-"
$ws.Range("D7").HorizontalAlignment = -4152      # xlRight
$ws.Range("D8").Value = "-"
$ws.Range("D8").HorizontalAlignment = -4152      # xlRight

# --- Restore the cursor position recorded in the saved view ---
$ws.Range("N5").Select() | Out-Null
